$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Helper: find a paragraph whose text starts with $needle and return it.
# ---------------------------------------------------------------------------
function Get-ParaByStart([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "Във вюто на листнатите ..." paragraph: split "Във вюто" so that "вюто"
#    is wrapped in proofErr spell-check markers (no formatting change).
# ---------------------------------------------------------------------------
$p1 = Get-ParaByStart("Във вюто на листнатите")
$xml1 = '<w:p ' + $wNs + '>' +
          '<w:pPr>' +
            '<w:pStyle w:val="ListParagraph"/>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
            '<w:contextualSpacing w:val="0"/>' +
          '</w:pPr>' +
          '<w:r><w:t xml:space="preserve">Във </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>вюто</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> на листнатите производителите, предлагам да добавим double touch върху някой който да те праща в edit формата за него. Това ще го направим след като имаме edit форма smile emoticon</w:t></w:r>' +
        '</w:p>'
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) "В детайлното вю ..." paragraph: keep the green highlight but drop the
#    _GoBack bookmark that used to sit at the start of it.
# ---------------------------------------------------------------------------
$p2 = Get-ParaByStart("В детайлното вю")
$xml2 = '<w:p ' + $wNs + '>' +
          '<w:pPr>' +
            '<w:pStyle w:val="ListParagraph"/>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
            '<w:contextualSpacing w:val="0"/>' +
            '<w:rPr><w:highlight w:val="green"/></w:rPr>' +
          '</w:pPr>' +
          '<w:r>' +
            '<w:rPr><w:highlight w:val="green"/></w:rPr>' +
            '<w:t>В детайлното вю трябва да добавим търсене на маршрут, когато се натисне адресът. Това ще го направя аз, след като си готова с другите неща по него.</w:t>' +
          '</w:r>' +
        '</w:p>'
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) "На Details page-а ..." paragraph: gains the green highlight and now
#    carries the _GoBack bookmark at the end of it (moved from paragraph 2).
# ---------------------------------------------------------------------------
$p3 = Get-ParaByStart("На Details page-а")
$xml3 = '<w:p ' + $wNs + '>' +
          '<w:pPr>' +
            '<w:pStyle w:val="ListParagraph"/>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
            '<w:contextualSpacing w:val="0"/>' +
            '<w:rPr><w:highlight w:val="green"/></w:rPr>' +
          '</w:pPr>' +
          '<w:r>' +
            '<w:rPr><w:highlight w:val="green"/></w:rPr>' +
            '<w:t>На Details page-а трябва да имаме бутон към Edit producer. Предлагам ти да го направиш, след като аз съм готов с Edit page-а.</w:t>' +
          '</w:r>' +
          '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
          '<w:bookmarkEnd w:id="0"/>' +
        '</w:p>'
$p3.Range.InsertXML($xml3)

Write-Host "Edit applied."
